$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.651.65"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").Value = "1.864.87"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.011"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.010"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4697"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3920"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  -3.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07990"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.997"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.859.06"
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.012"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.51%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06736"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001042"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").Value = "27.625.19"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.456"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.59%  "
$ws.Range("E24").Value = "  -1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.311"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "2.087.10"
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.143"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.433"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09484"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.617"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.298"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.335"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -8.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06057"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02228"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.327"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.196"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.009"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5950"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1883"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5646"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.920"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06759"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.111"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.44%  "
